# Update "Import"/"Export" monthly data (Data/monthly_data_total.xlsx) with
# the April-2019-onward refresh described in the commit message:
#   "Add files via upload / Update data from April 2019 for Import and Export"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Revised figures for the months that were already present (rows 161-166,
#    i.e. Apr-2019 .. Sep-2019). Dates (column A) are unchanged - only the
#    Import (B) and Export (C) values were refreshed.
# ---------------------------------------------------------------------------
$ws.Range("B161").Value = 42384.92
$ws.Range("C161").Value = 26025.68

$ws.Range("B162").Value = 46672.49
$ws.Range("C162").Value = 29847.36

$ws.Range("B163").Value = 41032.13
$ws.Range("C163").Value = 25017.07

$ws.Range("B164").Value = 40430.18
$ws.Range("C164").Value = 26227.1

$ws.Range("B165").Value = 39852.43
$ws.Range("C165").Value = 25981.42

$ws.Range("B166").Value = 37693.81
$ws.Range("C166").Value = 26007.97

# ---------------------------------------------------------------------------
# 2) Newly-populated months (rows 167-170, Oct-2019 .. Jan-2020). These cells
#    were previously blank, so the number formats have to be (re)applied
#    explicitly: column A uses the existing date format, B/C use the
#    existing 2-decimal number format, matching styles "1" and "2" already
#    used throughout the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A167").NumberFormat = "mmm-yy"
$ws.Range("A167").Value = 43739
$ws.Range("B167:C167").NumberFormat = "#,##0.00"
$ws.Range("B167").Value = 37241.53
$ws.Range("C167").Value = 26213.39

$ws.Range("A168").NumberFormat = "mmm-yy"
$ws.Range("A168").Value = 43770
$ws.Range("B168:C168").NumberFormat = "#,##0.00"
$ws.Range("B168").Value = 38101.94
$ws.Range("C168").Value = 25630.87

$ws.Range("A169").NumberFormat = "mmm-yy"
$ws.Range("A169").Value = 43800
$ws.Range("B169:C169").NumberFormat = "#,##0.00"
$ws.Range("B169").Value = 38577.34
$ws.Range("C169").Value = 27142.63

$ws.Range("A170").NumberFormat = "mmm-yy"
$ws.Range("A170").Value = 43831
$ws.Range("B170:C170").NumberFormat = "#,##0.00"
$ws.Range("B170").Value = 41146.86
$ws.Range("C170").Value = 25882.9

# Row 171 gains a (still empty) formatted cell in column B, mirroring the
# author's commit which pre-formatted the next row down the column.
$ws.Range("B171").NumberFormat = "#,##0.00"

# ---------------------------------------------------------------------------
# 3) View state: the author scrolled the sheet down and left the selection
#    on the newly-added C170 (last data entry) instead of the original C20.
# ---------------------------------------------------------------------------
$ws.Range("C170").Select()
$excel.ActiveWindow.ScrollRow = 149
$excel.ActiveWindow.ScrollColumn = 1
